# participante.xlsx update:
#  - add a new participant e-mail (row 4 / column D) with a mailto: hyperlink,
#    matching the style already used by the other e-mail cells
#  - move the active selection to D12
#  - set the print page setup (paper size / orientation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New participant e-mail for row 4 (currently empty D4), with hyperlink
$ws.Range("D4").Value = "gmysto6@gmail.com"
[void]$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:gmysto6@gmail.com")
$ws.Range("D4").Style = "Hyperlink"

# Update the saved selection/active cell
[void]$ws.Range("D12").Select()

# Page setup: A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
